$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 98
$ws.Cells.Item(2, 2).Value = '$(\langle pq \vert pq \rangle)_{3}$'
$ws.Cells.Item(2, 3).Value = 0.0008892826004847061
$ws.Cells.Item(3, 1).Value = 91
$ws.Cells.Item(3, 2).Value = '$(\langle pq \vert qp \rangle)_{2}$'
$ws.Cells.Item(3, 3).Value = 0.000678641705020092
$ws.Cells.Item(4, 1).Value = 0
$ws.Cells.Item(4, 2).Value = '(h$_{p}$)$_{0}$'
$ws.Cells.Item(4, 3).Value = 0.0004185367885651768
$ws.Cells.Item(5, 1).Value = 26
$ws.Cells.Item(5, 2).Value = 'type_3'
$ws.Cells.Item(5, 3).Value = 0.0003770353426672187
$ws.Cells.Item(6, 1).Value = 93
$ws.Cells.Item(6, 2).Value = '$(\langle rs \vert sr \rangle)_{2}$'
$ws.Cells.Item(6, 3).Value = 0.000318408844430717
$ws.Cells.Item(7, 1).Value = 12
$ws.Cells.Item(7, 2).Value = 'h$_{q}$'
$ws.Cells.Item(7, 3).Value = 0.000282804012915769
$ws.Cells.Item(8, 1).Value = 100
$ws.Cells.Item(8, 2).Value = '$(\langle rs\vert rs \rangle)_{3}$'
$ws.Cells.Item(8, 3).Value = 0.0002649591332523719
$ws.Cells.Item(9, 1).Value = 20
$ws.Cells.Item(9, 2).Value = '(h$_{rs}$)$_{2}$'
$ws.Cells.Item(9, 3).Value = 0.0001743158079944655
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = '(h$_{pq}$)$_{3}$'
$ws.Cells.Item(10, 3).Value = 0.0001603683158761672
$ws.Cells.Item(11, 1).Value = 78
$ws.Cells.Item(11, 2).Value = '$(\langle pq \vert rs \rangle)_{1}$'
$ws.Cells.Item(11, 3).Value = 0.0001187498436718664
$ws.Cells.Item(12, 1).Value = 97
$ws.Cells.Item(12, 2).Value = '$(\langle rr \vert rr \rangle)_{3}$'
$ws.Cells.Item(12, 3).Value = 0.0001179385892459306
$ws.Cells.Item(13, 1).Value = 85
$ws.Cells.Item(13, 2).Value = '$(\langle rs \vert sr \rangle)_{1}$'
$ws.Cells.Item(13, 3).Value = 0.000106636775606351
$ws.Cells.Item(14, 1).Value = 2
$ws.Cells.Item(14, 2).Value = '(h$_{p}$)$_{2}$'
$ws.Cells.Item(14, 3).Value = 0.00009139713873590455
$ws.Cells.Item(15, 1).Value = 3
$ws.Cells.Item(15, 2).Value = '(h$_{p}$)$_{3}$'
$ws.Cells.Item(15, 3).Value = 0.0000858107187398563
$ws.Cells.Item(16, 1).Value = 90
$ws.Cells.Item(16, 2).Value = '$(\langle pq \vert pq \rangle)_{2}$'
$ws.Cells.Item(16, 3).Value = 0.00007984438912437895
$ws.Cells.Item(17, 1).Value = 15
$ws.Cells.Item(17, 2).Value = '(h$_{r}$)$_{1}$'
$ws.Cells.Item(17, 3).Value = 0.0000675965346257002
$ws.Cells.Item(18, 1).Value = 67
$ws.Cells.Item(18, 2).Value = '$(\eta_{r})_{3}$'
$ws.Cells.Item(18, 3).Value = 0.0000651821297973602
$ws.Cells.Item(19, 1).Value = 92
$ws.Cells.Item(19, 2).Value = '$(\langle rs\vert rs \rangle)_{2}$'
$ws.Cells.Item(19, 3).Value = 0.00006089659296894452
$ws.Cells.Item(20, 1).Value = 21
$ws.Cells.Item(20, 2).Value = '(h$_{rs}$)$_{3}$'
$ws.Cells.Item(20, 3).Value = 0.00006082314601596453
$ws.Cells.Item(21, 1).Value = 32
$ws.Cells.Item(21, 2).Value = '$F_{q}^{\text{SCF}}$'
$ws.Cells.Item(21, 3).Value = 0.00006038503571868398
$ws.Cells.Item(22, 1).Value = 52
$ws.Cells.Item(22, 2).Value = '$(F_{p}^{\text{SCF}})_{2}$'
$ws.Cells.Item(22, 3).Value = 0.00004914754399614944
$ws.Cells.Item(23, 1).Value = 6
$ws.Cells.Item(23, 2).Value = '(h$_{pq}$)$_{2}$'
$ws.Cells.Item(23, 3).Value = 0.00004638017750896637
$ws.Cells.Item(24, 1).Value = 4
$ws.Cells.Item(24, 2).Value = '(h$_{pq}$)$_{0}$'
$ws.Cells.Item(24, 3).Value = 0.00004565585531285447
$ws.Cells.Item(25, 1).Value = 22
$ws.Cells.Item(25, 2).Value = 'h$_{s}$'
$ws.Cells.Item(25, 3).Value = 0.00004447190546208865
$ws.Cells.Item(26, 1).Value = 11
$ws.Cells.Item(26, 2).Value = '(h$_{pr}$)$_{3}$'
$ws.Cells.Item(26, 3).Value = 0.00004379944614525561
$ws.Cells.Item(27, 1).Value = 82
$ws.Cells.Item(27, 2).Value = '$(\langle pq \vert pq \rangle)_{1}$'
$ws.Cells.Item(27, 3).Value = 0.00003998616213711897
$ws.Cells.Item(28, 1).Value = 42
$ws.Cells.Item(28, 2).Value = '$F_{s}$'
$ws.Cells.Item(28, 3).Value = 0.00003995032543784376
$ws.Cells.Item(29, 1).Value = 36
$ws.Cells.Item(29, 2).Value = '$(F_{r}^{\text{SCF}})_{0}$'
$ws.Cells.Item(29, 3).Value = 0.00003852227837479863
$ws.Cells.Item(30, 1).Value = 17
$ws.Cells.Item(30, 2).Value = '(h$_{r}$)$_{3}$'
$ws.Cells.Item(30, 3).Value = 0.00003789972913217216
$ws.Cells.Item(31, 1).Value = 84
$ws.Cells.Item(31, 2).Value = '$(\langle rs\vert rs \rangle)_{1}$'
$ws.Cells.Item(31, 3).Value = 0.00003752580069846353
$ws.Cells.Item(32, 1).Value = 16
$ws.Cells.Item(32, 2).Value = '(h$_{r}$)$_{2}$'
$ws.Cells.Item(32, 3).Value = 0.00003750724468840098
$ws.Cells.Item(33, 1).Value = 40
$ws.Cells.Item(33, 2).Value = '$F_{s}^{\text{SCF}}$'
$ws.Cells.Item(33, 3).Value = 0.00003597268460202609
$ws.Cells.Item(34, 1).Value = 77
$ws.Cells.Item(34, 2).Value = '$(\langle rs \vert sr \rangle)_{0}$'
$ws.Cells.Item(34, 3).Value = 0.00003351675424879212
$ws.Cells.Item(35, 1).Value = 76
$ws.Cells.Item(35, 2).Value = '$(\langle rs\vert rs \rangle)_{0}$'
$ws.Cells.Item(35, 3).Value = 0.00002831506395982808
$ws.Cells.Item(36, 1).Value = 73
$ws.Cells.Item(36, 2).Value = '$\langle ss \vert ss \rangle$'
$ws.Cells.Item(36, 3).Value = 0.00002773139227179528
$ws.Cells.Item(37, 1).Value = 101
$ws.Cells.Item(37, 2).Value = '$(\langle rs \vert sr \rangle)_{3}$'
$ws.Cells.Item(37, 3).Value = 0.000027262452011274
$ws.Cells.Item(38, 1).Value = 71
$ws.Cells.Item(38, 2).Value = '$\langle qq \vert qq \rangle$'
$ws.Cells.Item(38, 3).Value = 0.00002695114288782441
$ws.Cells.Item(39, 1).Value = 74
$ws.Cells.Item(39, 2).Value = '$(\langle pq \vert pq \rangle)_{0}$'
$ws.Cells.Item(39, 3).Value = 0.00002643246565534268
$ws.Cells.Item(40, 1).Value = 99
$ws.Cells.Item(40, 2).Value = '$(\langle pq \vert qp \rangle)_{3}$'
$ws.Cells.Item(40, 3).Value = 0.0000229066411143425
$ws.Cells.Item(41, 1).Value = 83
$ws.Cells.Item(41, 2).Value = '$(\langle pq \vert qp \rangle)_{1}$'
$ws.Cells.Item(41, 3).Value = 0.00002123584403257723
$ws.Cells.Item(42, 1).Value = 86
$ws.Cells.Item(42, 2).Value = '$(\langle pq \vert rs \rangle)_{2}$'
$ws.Cells.Item(42, 3).Value = 0.00001975935653633256
$ws.Cells.Item(43, 1).Value = 1
$ws.Cells.Item(43, 2).Value = '(h$_{p}$)$_{1}$'
$ws.Cells.Item(43, 3).Value = 0.0000193747513273645
$ws.Cells.Item(44, 1).Value = 39
$ws.Cells.Item(44, 2).Value = '$(\eta_{r})_{0}$'
$ws.Cells.Item(44, 3).Value = 0.00001849147553331105
$ws.Cells.Item(45, 1).Value = 5
$ws.Cells.Item(45, 2).Value = '(h$_{pq}$)$_{1}$'
$ws.Cells.Item(45, 3).Value = 0.00001799840081848746
$ws.Cells.Item(46, 1).Value = 72
$ws.Cells.Item(46, 2).Value = '$(\langle rr \vert rr \rangle)_{0}$'
$ws.Cells.Item(46, 3).Value = 0.0000176564529973295
$ws.Cells.Item(47, 1).Value = 13
$ws.Cells.Item(47, 2).Value = 'h$_{qs}$'
$ws.Cells.Item(47, 3).Value = 0.00001683009207228131
$ws.Cells.Item(48, 1).Value = 10
$ws.Cells.Item(48, 2).Value = '(h$_{pr}$)$_{2}$'
$ws.Cells.Item(48, 3).Value = 0.00001679797074517646
$ws.Cells.Item(49, 1).Value = 34
$ws.Cells.Item(49, 2).Value = '$F_{q}$'
$ws.Cells.Item(49, 3).Value = 0.00001653591723972108
$ws.Cells.Item(50, 1).Value = 89
$ws.Cells.Item(50, 2).Value = '$(\langle rr \vert rr \rangle)_{2}$'
$ws.Cells.Item(50, 3).Value = 0.00001573152552394466
$ws.Cells.Item(51, 1).Value = 94
$ws.Cells.Item(51, 2).Value = '$(\langle pq \vert rs \rangle)_{3}$'
$ws.Cells.Item(51, 3).Value = 0.0000154443645710913
$ws.Cells.Item(52, 1).Value = 68
$ws.Cells.Item(52, 2).Value = '$(\langle pq \vert rs \rangle)_{0}$'
$ws.Cells.Item(52, 3).Value = 0.00001518898231305793
$ws.Cells.Item(53, 1).Value = 38
$ws.Cells.Item(53, 2).Value = '$(F_{r})_{0}$'
$ws.Cells.Item(53, 3).Value = 0.00001413317862715608
$ws.Cells.Item(54, 1).Value = 43
$ws.Cells.Item(54, 2).Value = '$\eta_{s}$'
$ws.Cells.Item(54, 3).Value = 0.00001296527535241432
$ws.Cells.Item(55, 1).Value = 29
$ws.Cells.Item(55, 2).Value = '$(\omega_{p})_{0}$'
$ws.Cells.Item(55, 3).Value = 0.00001247964623404556
$ws.Cells.Item(56, 1).Value = 80
$ws.Cells.Item(56, 2).Value = '$(\langle pp \vert pp \rangle)_{1}$'
$ws.Cells.Item(56, 3).Value = 0.000009489814348728603
$ws.Cells.Item(57, 1).Value = 75
$ws.Cells.Item(57, 2).Value = '$(\langle pq \vert qp \rangle)_{0}$'
$ws.Cells.Item(57, 3).Value = 0.00000900926744940274
$ws.Cells.Item(58, 1).Value = 9
$ws.Cells.Item(58, 2).Value = '(h$_{pr}$)$_{1}$'
$ws.Cells.Item(58, 3).Value = 0.000007919639178307435
$ws.Cells.Item(59, 1).Value = 50
$ws.Cells.Item(59, 2).Value = '$(F_{r})_{1}$'
$ws.Cells.Item(59, 3).Value = 0.000007365645153613365
$ws.Cells.Item(60, 1).Value = 60
$ws.Cells.Item(60, 2).Value = '$(F_{p}^{\text{SCF}})_{3}$'
$ws.Cells.Item(60, 3).Value = 0.000006898712550157791
$ws.Cells.Item(61, 1).Value = 70
$ws.Cells.Item(61, 2).Value = '$(\langle pp \vert pp \rangle)_{0}$'
$ws.Cells.Item(61, 3).Value = 0.000006611660775492922
$ws.Cells.Item(62, 1).Value = 8
$ws.Cells.Item(62, 2).Value = '(h$_{pr}$)$_{0}$'
$ws.Cells.Item(62, 3).Value = 0.000006487266848876169
$ws.Cells.Item(63, 1).Value = 44
$ws.Cells.Item(63, 2).Value = '$(F_{p}^{\text{SCF}})_{1}$'
$ws.Cells.Item(63, 3).Value = 0.000006204787670707869
$ws.Cells.Item(64, 1).Value = 14
$ws.Cells.Item(64, 2).Value = '(h$_{r}$)$_{0}$'
$ws.Cells.Item(64, 3).Value = 0.000006093062827095093
$ws.Cells.Item(65, 1).Value = 88
$ws.Cells.Item(65, 2).Value = '$(\langle pp \vert pp \rangle)_{2}$'
$ws.Cells.Item(65, 3).Value = 0.000005759952306854176
$ws.Cells.Item(66, 1).Value = 47
$ws.Cells.Item(66, 2).Value = '$(\eta_{p})_{1}$'
$ws.Cells.Item(66, 3).Value = 0.000005538180845518032
$ws.Cells.Item(67, 1).Value = 24
$ws.Cells.Item(67, 2).Value = 'type_1'
$ws.Cells.Item(67, 3).Value = 0.000004877284250885272
$ws.Cells.Item(68, 1).Value = 19
$ws.Cells.Item(68, 2).Value = '(h$_{rs}$)$_{1}$'
$ws.Cells.Item(68, 3).Value = 0.000004525621337882048
$ws.Cells.Item(69, 1).Value = 25
$ws.Cells.Item(69, 2).Value = 'type_2'
$ws.Cells.Item(69, 3).Value = 0.000004430655161964942
$ws.Cells.Item(70, 1).Value = 64
$ws.Cells.Item(70, 2).Value = '$(F_{r}^{\text{SCF}})_{3}$'
$ws.Cells.Item(70, 3).Value = 0.000004225423380743252
$ws.Cells.Item(71, 1).Value = 46
$ws.Cells.Item(71, 2).Value = '$(F_{p})_{1}$'
$ws.Cells.Item(71, 3).Value = 0.000004145679804425736
$ws.Cells.Item(72, 1).Value = 96
$ws.Cells.Item(72, 2).Value = '$(\langle pp \vert pp \rangle)_{3}$'
$ws.Cells.Item(72, 3).Value = 0.000003728544934711371
$ws.Cells.Item(73, 1).Value = 81
$ws.Cells.Item(73, 2).Value = '$(\langle rr \vert rr \rangle)_{1}$'
$ws.Cells.Item(73, 3).Value = 0.000003633500390251857
$ws.Cells.Item(74, 1).Value = 59
$ws.Cells.Item(74, 2).Value = '$(\eta_{r})_{2}$'
$ws.Cells.Item(74, 3).Value = 0.000003607939330915683
$ws.Cells.Item(75, 1).Value = 56
$ws.Cells.Item(75, 2).Value = '$(F_{r}^{\text{SCF}})_{2}$'
$ws.Cells.Item(75, 3).Value = 0.000003446142487066101
$ws.Cells.Item(76, 1).Value = 28
$ws.Cells.Item(76, 2).Value = '$(F_{p}^{\text{SCF}})_{0}$'
$ws.Cells.Item(76, 3).Value = 0.000003090075014405905
$ws.Cells.Item(77, 1).Value = 66
$ws.Cells.Item(77, 2).Value = '$(F_{r})_{3}$'
$ws.Cells.Item(77, 3).Value = 0.000002453736191715734
$ws.Cells.Item(78, 1).Value = 48
$ws.Cells.Item(78, 2).Value = '$(F_{r}^{\text{SCF}})_{1}$'
$ws.Cells.Item(78, 3).Value = 0.000002336712420874132
$ws.Cells.Item(79, 1).Value = 31
$ws.Cells.Item(79, 2).Value = '$(\eta_{p})_{0}$'
$ws.Cells.Item(79, 3).Value = 0.000002029162119557156
$ws.Cells.Item(80, 1).Value = 57
$ws.Cells.Item(80, 2).Value = '$(\omega_{r})_{2}$'
$ws.Cells.Item(80, 3).Value = 0.000001763645345810022
$ws.Cells.Item(81, 1).Value = 55
$ws.Cells.Item(81, 2).Value = '$(\eta_{p})_{2}$'
$ws.Cells.Item(81, 3).Value = 0.000001528031936514659
$ws.Cells.Item(82, 1).Value = 63
$ws.Cells.Item(82, 2).Value = '$(\eta_{p})_{3}$'
$ws.Cells.Item(82, 3).Value = 0.000001490608492816358
$ws.Cells.Item(83, 1).Value = 54
$ws.Cells.Item(83, 2).Value = '$(F_{p})_{2}$'
$ws.Cells.Item(83, 3).Value = 0.000001455622467277225
$ws.Cells.Item(84, 1).Value = 51
$ws.Cells.Item(84, 2).Value = '$(\eta_{r})_{1}$'
$ws.Cells.Item(84, 3).Value = 0.0000007030991053320424
$ws.Cells.Item(85, 1).Value = 58
$ws.Cells.Item(85, 2).Value = '$(F_{r})_{2}$'
$ws.Cells.Item(85, 3).Value = 0.0000003034298319814945
$ws.Cells.Item(86, 1).Value = 30
$ws.Cells.Item(86, 2).Value = '$(F_{p})_{0}$'
$ws.Cells.Item(86, 3).Value = 0.00000009003849366383107
$ws.Cells.Item(87, 1).Value = 79
$ws.Cells.Item(87, 2).Value = '$(\langle pq \vert sr \rangle)_{1}$'
$ws.Cells.Item(87, 3).Value = 0.00000004654496990812673
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = '$(\langle pq \vert sr \rangle)_{2}$'
$ws.Cells.Item(88, 3).Value = 0.00000003248439725662721
$ws.Cells.Item(89, 1).Value = 62
$ws.Cells.Item(89, 2).Value = '$(F_{p})_{3}$'
$ws.Cells.Item(89, 3).Value = 0.00000003181265347291359
$ws.Cells.Item(90, 1).Value = 95
$ws.Cells.Item(90, 2).Value = '$(\langle pq \vert sr \rangle)_{3}$'
$ws.Cells.Item(90, 3).Value = 0.0000000294778668239382
$ws.Cells.Item(91, 1).Value = 53
$ws.Cells.Item(91, 2).Value = '$(\omega_{p})_{2}$'
$ws.Cells.Item(91, 3).Value = 0.00000002063205841877757
$ws.Cells.Item(92, 1).Value = 49
$ws.Cells.Item(92, 2).Value = '$(\omega_{r})_{1}$'
$ws.Cells.Item(92, 3).Value = 0.00000002057912345873382
$ws.Cells.Item(93, 1).Value = 69
$ws.Cells.Item(93, 2).Value = '$(\langle pq \vert sr \rangle)_{0}$'
$ws.Cells.Item(93, 3).Value = 0.00000002037538446634959
$ws.Cells.Item(94, 1).Value = 23
$ws.Cells.Item(94, 2).Value = 'type_0'
$ws.Cells.Item(94, 3).Value = 0.00000001477521927925031
$ws.Cells.Item(95, 1).Value = 18
$ws.Cells.Item(95, 2).Value = '(h$_{rs}$)$_{0}$'
$ws.Cells.Item(95, 3).Value = 0.00000001454319361954811
$ws.Cells.Item(96, 1).Value = 37
$ws.Cells.Item(96, 2).Value = '$(\omega_{r})_{0}$'
$ws.Cells.Item(96, 3).Value = 0.00000001265740218356881
$ws.Cells.Item(97, 1).Value = 65
$ws.Cells.Item(97, 2).Value = '$(\omega_{r})_{3}$'
$ws.Cells.Item(97, 3).Value = 0.00000001244621199806713
$ws.Cells.Item(98, 1).Value = 27
$ws.Cells.Item(98, 2).Value = '$\mathbf{b}$'
$ws.Cells.Item(98, 3).Value = 0.000000008700893674223822
$ws.Cells.Item(99, 1).Value = 45
$ws.Cells.Item(99, 2).Value = '$(\omega_{p})_{1}$'
$ws.Cells.Item(99, 3).Value = 0.000000005005227185111558
$ws.Cells.Item(100, 1).Value = 33
$ws.Cells.Item(100, 2).Value = '$\omega_{q}$'
$ws.Cells.Item(100, 3).Value = 0.000000003500057336255808
$ws.Cells.Item(101, 1).Value = 41
$ws.Cells.Item(101, 2).Value = '$\omega_{s}$'
$ws.Cells.Item(101, 3).Value = 0.000000002639809626868773
$ws.Cells.Item(102, 1).Value = 61
$ws.Cells.Item(102, 2).Value = '$(\omega_{p})_{3}$'
$ws.Cells.Item(102, 3).Value = 0.000000001644199942163898
$ws.Cells.Item(103, 1).Value = 35
$ws.Cells.Item(103, 2).Value = '$\eta_{q}$'
$ws.Cells.Item(103, 3).Value = 0.0000000004317328074389671
